# B6-PowerPoint.pptx edit — "Mon, May 04, 2020  4:08:58 PM"
#
# Two logical changes are reproduced here:
#
# 1. The table on slides 14, 15 and 16 is re-styled from the deck's
#    default table style ({DE7B04B8-A9E9-4F7E-8932-AB00B0E485F6}) to the
#    built-in style {14ED5A61-1697-4B48-B9E9-18E81FFAD7D8}
#    (PowerPoint's "Medium Style 2 - Accent 1").
#
# 2. The presentation's theme palette is swapped from the custom
#    "Integral" / "Red Violet" colour scheme to the stock Office theme
#    palette (the font scheme and the format scheme were already shared
#    between the two theme parts, so only the 12 theme colours differ).

$p = $ppt.ActivePresentation

# --- 1. Re-apply the table style on the three slides that contain a table ---
$tableSlides = @(14, 15, 16)
foreach ($slideIndex in $tableSlides) {
    $slide = $p.Slides.Item($slideIndex)
    $tableShape = $slide.Shapes.Item(1)
    if ($tableShape.HasTable) {
        $tableShape.Table.ApplyStyle("{14ED5A61-1697-4B48-B9E9-18E81FFAD7D8}")
    }
}

# --- 2. Swap the theme colour scheme to the stock "Office" palette ---
# MsoThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# RGB values are passed as the usual COM 0xBBGGRR integers.
$slide1 = $p.Slides.Item(1)
$colorScheme = $slide1.ThemeColorScheme
$colorScheme.Item(1).RGB = 0           # dk1      000000
$colorScheme.Item(2).RGB = 16777215    # lt1      FFFFFF
$colorScheme.Item(3).RGB = 6968388     # dk2      44546A
$colorScheme.Item(4).RGB = 15132391    # lt2      E7E6E6
$colorScheme.Item(5).RGB = 13998939    # accent1  5B9BD5
$colorScheme.Item(6).RGB = 3243501     # accent2  ED7D31
$colorScheme.Item(7).RGB = 10855845    # accent3  A5A5A5
$colorScheme.Item(8).RGB = 49407       # accent4  FFC000
$colorScheme.Item(9).RGB = 12874308    # accent5  4472C4
$colorScheme.Item(10).RGB = 4697456    # accent6  70AD47
$colorScheme.Item(11).RGB = 12673797   # hlink    0563C1
$colorScheme.Item(12).RGB = 7491477    # folHlink 954F72
